$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Add new header cell J1 with the same formatting as I1 (bold, centered, bordered)
$ws.Range("J1").Value = "23-jun"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Fill in the new "23-jun" price column (J2:J25)
$values = @(41.21, 23.6, 15.66, 12.29, 12.2, 12.51, 19.63, 26, 38.41, 14.96, 0.02, -0.01, -1.98, -8.890000000000001, -17.32, -9.529999999999999, -0.06, 0, 20.98, 47.44, 94.98, 100.15, 106.5, 94.55)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $values[$i]
}
